$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.039.48"
$ws.Range("E2").Value = "  -0.33%  "

$ws.Range("D3").Value = "1.828.73"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("D5").Value = "240.81"

$ws.Range("D6").Value = "0.6226"
$ws.Range("E6").Value = "  -6.06%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "44.51"
$ws.Range("E8").Value = "  +6.10%  "

$ws.Range("D9").Value = "0.07415"
$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").Value = "0.2919"
$ws.Range("E10").Value = "  -0.41%  "

$ws.Range("E11").Value = "  +0.03%  "

$ws.Range("D12").Value = "0.07603"
$ws.Range("E12").Value = "  -1.73%  "

$ws.Range("D13").Value = "1.831.27"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("E14").Value = "  -0.38%  "

$ws.Range("D15").Value = "0.6624"
$ws.Range("E15").Value = "  -0.68%  "

$ws.Range("D16").Value = "82.03"
$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("D17").Value = "0.000009143"
$ws.Range("E17").Value = "  +9.71%  "

$ws.Range("D18").Value = "6.009"
$ws.Range("E18").Value = "  -1.16%  "

$ws.Range("D19").Value = "29.043.83"
$ws.Range("E19").Value = "  -0.39%  "

$ws.Range("D20").Value = "2.079.96"
$ws.Range("E20").Value = "  -0.75%  "

$ws.Range("D21").Value = "225.13"
$ws.Range("E21").Value = "  -0.76%  "

$ws.Range("E22").Value = "  -0.83%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").Value = "7.174"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").Value = "159.31"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").Value = "8.411"
$ws.Range("E27").Value = "  -2.12%  "

$ws.Range("D28").Value = "0.1356"
$ws.Range("E28").Value = "  -3.02%  "

$ws.Range("D29").Value = "17.79"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").Value = "1.496"
$ws.Range("E30").Value = "  -0.86%  "

$ws.Range("E31").Value = "  -1.17%  "

$ws.Range("D32").Value = "4.025"
$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("E33").Value = "  +0.81%  "

$ws.Range("D34").Value = "0.05238"
$ws.Range("E34").Value = "  -1.29%  "

$ws.Range("D35").Value = "1.835"
$ws.Range("E35").Value = "  -1.68%  "

$ws.Range("D36").Value = "0.7346"
$ws.Range("E36").Value = "  -2.09%  "

$ws.Range("E37").Value = "  +1.40%  "

$ws.Range("E38").Value = "  +0.28%  "

$ws.Range("D39").Value = "1.280.41"
$ws.Range("E39").Value = "  +0.23%  "

$ws.Range("E40").Value = "  +0.81%  "

$ws.Range("D41").Value = "0.01781"
$ws.Range("E41").Value = "  -0.65%  "

$ws.Range("D42").Value = "6.320"
$ws.Range("E42").Value = "  +7.11%  "

$ws.Range("D43").Value = "0.8944"
$ws.Range("E43").Value = "  -3.48%  "

$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").Value = "101.64"
$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("D46").Value = "1.977.88"
$ws.Range("E46").Value = "  -0.87%  "

$ws.Range("E47").Value = "  -0.60%  "

$ws.Range("D48").Value = "63.63"
$ws.Range("E48").Value = "  +0.87%  "

$ws.Range("D49").Value = "0.00000000119"
$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("E50").Value = "  -3.14%  "

$ws.Range("D51").Value = "0.3963"
$ws.Range("E51").Value = "  -1.27%  "
